$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
$ws.Activate()

# Update the F:I result columns for rows 14-60 with the refreshed 18/04/2023 run data
# (replacing the previous 11/04/2023-13/04/2023 run captured before the corrected commit)
$ws.Range('F14').Value = '18/04/2023'
$ws.Range('G14').Value = '18/04/2023 02:24:56'
$ws.Range('H14').Value = '0fc88a0b96a2bfd4'
$ws.Range('I14').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.1f6d5429a6^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F15').Value = '18/04/2023'
$ws.Range('G15').Value = '18/04/2023 02:24:59'
$ws.Range('H15').Value = 'b1e83066d95c4b7d'
$ws.Range('I15').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.e2000b05b5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F16').Value = '18/04/2023'
$ws.Range('G16').Value = '18/04/2023 02:25:01'
$ws.Range('H16').Value = 'b5b327df15660c4a'
$ws.Range('I16').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.2eff41cf5b^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F17').Value = '18/04/2023'
$ws.Range('G17').Value = '18/04/2023 02:25:04'
$ws.Range('H17').Value = '759dfbf90b48fb62'
$ws.Range('I17').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.e254e75390^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G18').Value = '18/04/2023 02:27:07'
$ws.Range('H18').Value = 'cd73a73a3f86ff19'
$ws.Range('G19').Value = '18/04/2023 02:26:24'
$ws.Range('H19').Value = '6691efc8f329a433'
$ws.Range('F20').Value = '18/04/2023'
$ws.Range('G20').Value = '18/04/2023 02:27:13'
$ws.Range('H20').Value = '0fb717dd03b12646'
$ws.Range('F21').Value = '18/04/2023'
$ws.Range('G21').Value = '18/04/2023 02:26:26'
$ws.Range('H21').Value = '1be4693140518dc2'
$ws.Range('G24').Value = '18/04/2023 02:26:40'
$ws.Range('H24').Value = '995b3a4048f56911'
$ws.Range('I24').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.5f99f0dbb2^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G26').Value = '18/04/2023 02:26:42'
$ws.Range('H26').Value = '3faaf97bb47479c0'
$ws.Range('I26').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.648844608b^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G27').Value = '18/04/2023 02:26:45'
$ws.Range('H27').Value = '94ce2ac31869c8c3'
$ws.Range('I27').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.71b06f50d6^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G28').Value = '18/04/2023 02:26:47'
$ws.Range('H28').Value = 'cdd8c96a0e0a1a1a'
$ws.Range('I28').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.2658463ced^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G29').Value = '18/04/2023 02:26:51'
$ws.Range('H29').Value = '0ea713d0056944c9'
$ws.Range('I29').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.059c911179^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G30').Value = '18/04/2023 02:26:53'
$ws.Range('H30').Value = '4e341843397d917b'
$ws.Range('I30').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.c74fd9972c^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G31').Value = '18/04/2023 02:26:56'
$ws.Range('H31').Value = 'd359f91842537fa1'
$ws.Range('I31').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.0fa2797dce^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G32').Value = '18/04/2023 02:26:58'
$ws.Range('H32').Value = '3c2e18953ca1039f'
$ws.Range('I32').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.5d7f71f85a^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G33').Value = '18/04/2023 02:27:01'
$ws.Range('H33').Value = '2f143c8db4714910'
$ws.Range('I33').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.10ed52b182^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('G34').Value = '18/04/2023 02:27:04'
$ws.Range('H34').Value = 'd450358b15f42d5a'
$ws.Range('I34').Value = '2.16.840.1.113883.2.9.2.150.4.4.e883906759350f874e146c909b921a54c8d82a646997feb46ad26be1982458a2.f01e5d4b9e^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F36').Value = '18/04/2023'
$ws.Range('G36').Value = '18/04/2023 02:25:07'
$ws.Range('H36').Value = '45aa58abe73d74ce'
$ws.Range('I36').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.40e5516f77^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F38').Value = '18/04/2023'
$ws.Range('G38').Value = '18/04/2023 02:25:13'
$ws.Range('H38').Value = '70a2b35b2be236af'
$ws.Range('I38').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.b1e8faaae5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F39').Value = '18/04/2023'
$ws.Range('G39').Value = '18/04/2023 02:25:16'
$ws.Range('H39').Value = 'aea426a07fe0c28d'
$ws.Range('I39').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.8276916eb9^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F40').Value = '18/04/2023'
$ws.Range('G40').Value = '18/04/2023 02:25:18'
$ws.Range('H40').Value = '5c1afb00d3b5c0b5'
$ws.Range('I40').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.e05cc9c04c^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F41').Value = '18/04/2023'
$ws.Range('G41').Value = '18/04/2023 02:25:21'
$ws.Range('H41').Value = '433d3769d30fc808'
$ws.Range('I41').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.b7e446d242^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F42').Value = '18/04/2023'
$ws.Range('G42').Value = '18/04/2023 02:25:24'
$ws.Range('H42').Value = '01f4f6f24f46ca1a'
$ws.Range('I42').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.d59af4f6a1^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F43').Value = '18/04/2023'
$ws.Range('G43').Value = '18/04/2023 02:25:26'
$ws.Range('H43').Value = '0f7d2203e9090120'
$ws.Range('I43').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.b66d83b0d4^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F44').Value = '18/04/2023'
$ws.Range('G44').Value = '18/04/2023 02:25:29'
$ws.Range('H44').Value = '69a464ce643392e3'
$ws.Range('I44').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.e05882c107^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F45').Value = '18/04/2023'
$ws.Range('G45').Value = '18/04/2023 02:25:32'
$ws.Range('H45').Value = '1a23e82ed002466b'
$ws.Range('I45').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.ae996f4595^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F46').Value = '18/04/2023'
$ws.Range('G46').Value = '18/04/2023 02:25:39'
$ws.Range('H46').Value = '6e4ac6b273e93f59'
$ws.Range('I46').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.1f4c45ddea^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F47').Value = '18/04/2023'
$ws.Range('G47').Value = '18/04/2023 02:25:42'
$ws.Range('H47').Value = 'c359782383da95ad'
$ws.Range('I47').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.699b4c7aea^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F48').Value = '18/04/2023'
$ws.Range('G48').Value = '18/04/2023 02:25:45'
$ws.Range('H48').Value = '169000e8ca8a91ba'
$ws.Range('I48').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.2b8247e090^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F49').Value = '18/04/2023'
$ws.Range('G49').Value = '18/04/2023 02:25:48'
$ws.Range('H49').Value = '6f870121cd403f8b'
$ws.Range('I49').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.ef15b98aa5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F50').Value = '18/04/2023'
$ws.Range('G50').Value = '18/04/2023 02:25:50'
$ws.Range('H50').Value = 'e4c4451419020fdf'
$ws.Range('I50').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.63442d98a9^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F51').Value = '18/04/2023'
$ws.Range('G51').Value = '18/04/2023 02:25:53'
$ws.Range('H51').Value = '923f2d7bffc67a2b'
$ws.Range('I51').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.c578f27813^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F52').Value = '18/04/2023'
$ws.Range('G52').Value = '18/04/2023 02:25:55'
$ws.Range('H52').Value = '20caa40a5ba71c78'
$ws.Range('I52').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.96e8b6cf46^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F53').Value = '18/04/2023'
$ws.Range('G53').Value = '18/04/2023 02:25:59'
$ws.Range('H53').Value = 'a00c58f178e30daf'
$ws.Range('I53').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.5a2fc508ce^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F54').Value = '18/04/2023'
$ws.Range('G54').Value = '18/04/2023 02:26:03'
$ws.Range('H54').Value = '0377e7e31dd94f77'
$ws.Range('I54').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.94b1186b29^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F55').Value = '18/04/2023'
$ws.Range('G55').Value = '18/04/2023 02:26:08'
$ws.Range('H55').Value = '8b5585d3d4985efd'
$ws.Range('I55').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.aeab20db06^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F56').Value = '18/04/2023'
$ws.Range('G56').Value = '18/04/2023 02:26:11'
$ws.Range('H56').Value = '0ed3c7a6f8d4c4a5'
$ws.Range('I56').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.62d9d4ca13^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F57').Value = '18/04/2023'
$ws.Range('G57').Value = '18/04/2023 02:26:14'
$ws.Range('H57').Value = 'd11a11feb9af0ec7'
$ws.Range('I57').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.a6b4ba5ce1^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F58').Value = '18/04/2023'
$ws.Range('G58').Value = '18/04/2023 02:26:17'
$ws.Range('H58').Value = '5de2a7fc736315f9'
$ws.Range('I58').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.904d5b74d2^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F59').Value = '18/04/2023'
$ws.Range('G59').Value = '18/04/2023 02:26:19'
$ws.Range('H59').Value = '7da4cf147d0ffc05'
$ws.Range('I59').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.f0517e623b^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'
$ws.Range('F60').Value = '18/04/2023'
$ws.Range('G60').Value = '18/04/2023 02:26:22'
$ws.Range('H60').Value = '6161b033763ec183'
$ws.Range('I60').Value = '2.16.840.1.113883.2.9.2.150.4.4.05749625b2397823fb6e622bd34c86f95d4204a8d88eda08af3467630d91b2bd.32b6742cb5^^^^urn:ihe:iti:xdw:2013:workflowInstanceId'

# Restore the reviewer's last selection on the refreshed range
$ws.Range("F55:I60").Select()
